{"js": "// Apply text replacements for updated date + equations.\nconst pairs = [\n  [\"2025-02-27 Thursday\", \"2025-02-28 Friday\"],\n  [\"297\u00d75=1485\", \"933\u00d73=2799\"],\n  [\"869\u00d72=1738\", \"234\u00d79=2106\"],\n  [\"502\u00d76=3012\", \"137\u00d73=411\"],\n  [\"362\u00d77=2534\", \"335\u00d77=2345\"],\n  [\"160\u00d79=1440\", \"479\u00d75=2395\"],\n  [\"586\u00d79=5274\", \"550\u00d79=4950\"],\n  [\"339\u00d72=678\", \"375\u00d76=2250\"],\n  [\"816\u00d77=5712\", \"680\u00d75=3400\"],\n  [\"591\u00d77=4137\", \"275\u00d74=1100\"],\n  [\"394\u00d77=2758\", \"694\u00d77=4858\"],\n  [\"253\u00d73=759\", \"883\u00d76=5298\"],\n  [\"928\u00d78=7424\", \"476\u00d79=4284\"],\n  [\"571\u00d74=2284\", \"716\u00d74=2864\"],\n  [\"116\u00d78=928\", \"794\u00d75=3970\"],\n  [\"212\u00d73=636\", \"428\u00d77=2996\"],\n  [\"608\u00d79=5472\", \"843\u00d79=7587\"],\n  [\"269\u00d72=538\", \"649\u00d78=5192\"],\n  [\"254\u00d76=1524\", \"120\u00d76=720\"],\n  [\"738\u00d76=4428\", \"514\u00d78=4112\"],\n  [\"535\u00d76=3210\", \"564\u00d72=1128\"],\n  [\"644\u00d77=4508\", \"424\u00d78=3392\"],\n  [\"259\u00d79=2331\", \"951\u00d75=4755\"],\n  [\"686\u00d76=4116\", \"423\u00d72=846\"],\n  [\"247\u00d79=2223\", \"144\u00d77=1008\"],\n  [\"767\u00d77=5369\", \"692\u00d72=1384\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Update the worksheet date and the 25 multiplication equations.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-27 Thursday\", \"2025-02-28 Friday\"),\n    @(\"297\u00d75=1485\", \"933\u00d73=2799\"),\n    @(\"869\u00d72=1738\", \"234\u00d79=2106\"),\n    @(\"502\u00d76=3012\", \"137\u00d73=411\"),\n    @(\"362\u00d77=2534\", \"335\u00d77=2345\"),\n    @(\"160\u00d79=1440\", \"479\u00d75=2395\"),\n    @(\"586\u00d79=5274\", \"550\u00d79=4950\"),\n    @(\"339\u00d72=678\", \"375\u00d76=2250\"),\n    @(\"816\u00d77=5712\", \"680\u00d75=3400\"),\n    @(\"591\u00d77=4137\", \"275\u00d74=1100\"),\n    @(\"394\u00d77=2758\", \"694\u00d77=4858\"),\n    @(\"253\u00d73=759\", \"883\u00d76=5298\"),\n    @(\"928\u00d78=7424\", \"476\u00d79=4284\"),\n    @(\"571\u00d74=2284\", \"716\u00d74=2864\"),\n    @(\"116\u00d78=928\", \"794\u00d75=3970\"),\n    @(\"212\u00d73=636\", \"428\u00d77=2996\"),\n    @(\"608\u00d79=5472\", \"843\u00d79=7587\"),\n    @(\"269\u00d72=538\", \"649\u00d78=5192\"),\n    @(\"254\u00d76=1524\", \"120\u00d76=720\"),\n    @(\"738\u00d76=4428\", \"514\u00d78=4112\"),\n    @(\"535\u00d76=3210\", \"564\u00d72=1128\"),\n    @(\"644\u00d77=4508\", \"424\u00d78=3392\"),\n    @(\"259\u00d79=2331\", \"951\u00d75=4755\"),\n    @(\"686\u00d76=4116\", \"423\u00d72=846\"),\n    @(\"247\u00d79=2223\", \"144\u00d77=1008\"),\n    @(\"767\u00d77=5369\", \"692\u00d72=1384\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
